$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-12
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05)
$ws.Range("C2:C12").Value = 45174
